# Updates market/profit figures (currentAveragePrice, LevePrice, LeveProfit columns)
# across several Leve sheets, per scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H87").Value = 41152.668
$ws.Range("J87").Value = 41152.668
$ws.Range("L87").Value = 41152.668
$ws.Range("N87").Value = -43648.668

$ws.Range("H90").Value = 41152.668
$ws.Range("J90").Value = 41152.668
$ws.Range("L90").Value = 123458.004
$ws.Range("N90").Value = -135938.004

$ws.Range("H113").Value = 76928660
$ws.Range("I113").Value = 111115390
$ws.Range("K113").Value = 111115390
$ws.Range("M113").Value = -111112136

$ws.Range("H129").Value = 1057.4468
$ws.Range("J129").Value = 1219.4736
$ws.Range("L129").Value = 3658.4208
$ws.Range("N129").Value = -13658.4208

$ws.Range("H137").Value = 1966.75
$ws.Range("I137").Value = 1733.1111
$ws.Range("K137").Value = 5199.3333
$ws.Range("M137").Value = -2649.3333

$ws.Range("H138").Value = 34486164
$ws.Range("J138").Value = 4719.7144
$ws.Range("L138").Value = 14159.1432
$ws.Range("N138").Value = -24439.1432

# --- ARM sheet ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2226.6155
$ws.Range("I61").Value = 1654.1818
$ws.Range("K61").Value = 1654.1818
$ws.Range("M61").Value = -1442.1818

$ws.Range("H74").Value = 83335850
$ws.Range("I74").Value = 111112250
$ws.Range("J74").Value = 6666.6665
$ws.Range("K74").Value = 111112250
$ws.Range("L74").Value = 6666.6665
$ws.Range("M74").Value = -111111376
$ws.Range("N74").Value = -8414.666499999999

$ws.Range("H77").Value = 83335850
$ws.Range("I77").Value = 111112250
$ws.Range("J77").Value = 6666.6665
$ws.Range("K77").Value = 555561250
$ws.Range("L77").Value = 33333.3325
$ws.Range("M77").Value = -555556882
$ws.Range("N77").Value = -42069.3325

$ws.Range("H86").Value = 48000
$ws.Range("J86").Value = 48000
$ws.Range("L86").Value = 48000
$ws.Range("N86").Value = -50372

$ws.Range("H89").Value = 48000
$ws.Range("J89").Value = 48000
$ws.Range("L89").Value = 144000
$ws.Range("N89").Value = -155856

$ws.Range("H102").Value = 648.4761999999999
$ws.Range("I102").Value = 605.95
$ws.Range("K102").Value = 605.95
$ws.Range("M102").Value = 1016.05

$ws.Range("H132").Value = 15588.917
$ws.Range("I132").Value = 1338
$ws.Range("J132").Value = 47977.363
$ws.Range("K132").Value = 4014
$ws.Range("L132").Value = 143932.089
$ws.Range("M132").Value = -1484
$ws.Range("N132").Value = -148992.089

$ws.Range("H136").Value = 2226.6155
$ws.Range("I136").Value = 1654.1818
$ws.Range("K136").Value = 4962.5454
$ws.Range("M136").Value = -2412.5454

# --- BSM sheet ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H124").Value = 40693.332
$ws.Range("J124").Value = 40693.332
$ws.Range("L124").Value = 40693.332
$ws.Range("N124").Value = -50513.332

# --- CRP sheet ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H80").Value = 27128
$ws.Range("J80").Value = 27128
$ws.Range("L80").Value = 27128
$ws.Range("N80").Value = -29374

$ws.Range("H83").Value = 27128
$ws.Range("J83").Value = 27128
$ws.Range("L83").Value = 81384
$ws.Range("N83").Value = -92616

$ws.Range("H99").Value = 23812942
$ws.Range("J99").Value = 83338080
$ws.Range("L99").Value = 83338080
$ws.Range("N99").Value = -83341076

$ws.Range("H105").Value = 13890016
$ws.Range("I105").Value = 31250788
$ws.Range("J105").Value = 1398.8
$ws.Range("K105").Value = 31250788
$ws.Range("L105").Value = 1398.8
$ws.Range("M105").Value = -31249041
$ws.Range("N105").Value = -4892.8

$ws.Range("H107").Value = 1194.9678
$ws.Range("I107").Value = 928.46155
$ws.Range("J107").Value = 1387.4445
$ws.Range("K107").Value = 928.46155
$ws.Range("L107").Value = 1387.4445
$ws.Range("M107").Value = 991.53845
$ws.Range("N107").Value = -5227.4445

$ws.Range("H122").Value = 1483.6
$ws.Range("I122").Value = 1390.3334
$ws.Range("J122").Value = 1623.5
$ws.Range("K122").Value = 4171.0002
$ws.Range("L122").Value = 4870.5
$ws.Range("M122").Value = -1721.0002
$ws.Range("N122").Value = -9770.5

$ws.Range("H125").Value = 27000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 27000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 27000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -31920

$ws.Range("H126").Value = 23812942
$ws.Range("J126").Value = 83338080
$ws.Range("L126").Value = 250014240
$ws.Range("N126").Value = -250019180

$ws.Range("H132").Value = 3321
$ws.Range("I132").Value = 2605.25
$ws.Range("J132").Value = 5229.6665
$ws.Range("K132").Value = 7815.75
$ws.Range("L132").Value = 15688.9995
$ws.Range("M132").Value = -5285.75
$ws.Range("N132").Value = -20748.9995

$ws.Range("H134").Value = 992.4286
$ws.Range("I134").Value = 649.3077
$ws.Range("J134").Value = 1550
$ws.Range("K134").Value = 1947.9231
$ws.Range("L134").Value = 4650
$ws.Range("M134").Value = 587.0769
$ws.Range("N134").Value = -9720

# --- CUL sheet ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 741.52
$ws.Range("J131").Value = 772.413
$ws.Range("L131").Value = 2317.239
$ws.Range("N131").Value = -12397.239

# --- GSM sheet ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H15").Value = 21000
$ws.Range("J15").Value = 21000
$ws.Range("L15").Value = 21000
$ws.Range("N15").Value = -21576

$ws.Range("H81").Value = 21000
$ws.Range("J81").Value = 21000
$ws.Range("L81").Value = 21000
$ws.Range("N81").Value = -22996

$ws.Range("H84").Value = 21000
$ws.Range("J84").Value = 21000
$ws.Range("L84").Value = 63000
$ws.Range("N84").Value = -72984

$ws.Range("H113").Value = 2493.8235
$ws.Range("I113").Value = 2018.7778
$ws.Range("K113").Value = 2018.7778
$ws.Range("M113").Value = 151.2221999999999

$ws.Range("H126").Value = 4517.353
$ws.Range("I126").Value = 3925
$ws.Range("K126").Value = 11775
$ws.Range("M126").Value = -9305

# --- LTW sheet ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 4762.857
$ws.Range("I61").Value = 1720.5555
$ws.Range("K61").Value = 1720.5555
$ws.Range("M61").Value = -1518.5555

$ws.Range("H74").Value = 38400
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 38400
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H80").Value = 48000
$ws.Range("J80").Value = 48000
$ws.Range("L80").Value = 48000
$ws.Range("N80").Value = -50246

$ws.Range("H83").Value = 48000
$ws.Range("J83").Value = 48000
$ws.Range("L83").Value = 144000
$ws.Range("N83").Value = -155232

$ws.Range("H113").Value = 4762.857
$ws.Range("I113").Value = 1720.5555
$ws.Range("K113").Value = 1720.5555
$ws.Range("M113").Value = 449.4445000000001

$ws.Range("H124").Value = 30214.5
$ws.Range("J124").Value = 30214.5
$ws.Range("L124").Value = 30214.5
$ws.Range("N124").Value = -40034.5

$ws.Range("H125").Value = 28750
$ws.Range("J125").Value = 28750
$ws.Range("L125").Value = 28750
$ws.Range("N125").Value = -38590

$ws.Range("H132").Value = 671763.4
$ws.Range("I132").Value = 1340693.6
$ws.Range("K132").Value = 4022080.8
$ws.Range("M132").Value = -4019550.8

# --- WVR sheet ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H20").Value = 4000
$ws.Range("J20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("N20").Value = -4480

$ws.Range("H39").Value = 7285.7144
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 7666.6665
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 7666.6665
$ws.Range("M39").Value = -4587
$ws.Range("N39").Value = -8492.666499999999

$ws.Range("H75").Value = 25500
$ws.Range("J75").Value = 25500
$ws.Range("L75").Value = 25500
$ws.Range("N75").Value = -27372

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H78").Value = 25500
$ws.Range("J78").Value = 25500
$ws.Range("L78").Value = 76500
$ws.Range("N78").Value = -85860

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H113").Value = 9009186
$ws.Range("I113").Value = 250
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
